# Generate Report for Handback
# Updates the zh-cn and de-de localization-status sheets: the row for the
# 64bd9093-... file now failed handback transform (instead of being ready
# for handoff), record the error detail explaining the mismatched file
# name, and widen the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

$zhError = "Handback file name: qqwkdyb0.wag is different with handoff file name: 64bd9093-5465-443c-9d66-ce076fce8778.2dcff4d7a67333b2c86b16a494994cf243f3ab8c.zh-cn."
$deError  = "Handback file name: qqwkdyb0.wag is different with handoff file name: 64bd9093-5465-443c-9d66-ce076fce8778.2dcff4d7a67333b2c86b16a494994cf243f3ab8c.de-de."

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("P3").Value = $zhError
$wsZh.Columns.Item(16).ColumnWidth = 39.17

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("P3").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = 39.17

# The Overview sheet's per-locale status columns mirror the same "Status"
# text via the shared string table, so they need to move in lockstep.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
